$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "EasyShellTest().check_main_window(True)"
$ws.Range("B10").Value = "EasyShellTest().check_main_window(False)"
$ws.Range("C10").Clear()
$ws.Range("B10").Select()
